# TARUGO YESO DISMAY price list - fix values per "exceeded request" bug fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a hard-coded date serial (formerly 45310 / 2024-01-19) -> bump one day.
$ws.Range("A1").Value = 45311

# Price corrections for the two tarugo bag sizes.
$ws.Range("D28").Value = 25405.6
$ws.Range("D29").Value = 16946.6
